$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='64.295.22'; E='  +2.07%  '},
    @{Row=3; D='2.631.48'; E='  +0.35%  '},
    @{Row=4; E='  -0.01%  '},
    @{Row=5; D='603.77'; E='  +0.26%  '},
    @{Row=6; E='  +3.83%  '},
    @{Row=7; E='  -0.06%  '},
    @{Row=8; E='  +1.12%  '},
    @{Row=9; E='  +2.57%  '},
    @{Row=10; E='  +3.32%  '},
    @{Row=11; D='0.387'; E='  +7.10%  '},
    @{Row=12; E='  -0.64%  '},
    @{Row=13; D='27.67'; E='  +2.18%  '},
    @{Row=14; D='3.103.38'; E='  +0.34%  '},
    @{Row=15; D='64.161.36'; E='  +2.07%  '},
    @{Row=16; E='  +4.23%  '},
    @{Row=17; D='2.631.25'; E='  +0.22%  '},
    @{Row=18; E='  +8.30%  '},
    @{Row=19; D='4.66'; E='  +4.59%  '},
    @{Row=20; D='349.60'; E='  +3.23%  '},
    @{Row=21; D='6.99'; E='  +2.31%  '},
    @{Row=22; E='  +0.02%  '},
    @{Row=23; E='  +3.20%  '},
    @{Row=24; D='66.73'; E='  +0.53%  '},
    @{Row=25; E='  +15.92%  '},
    @{Row=26; D='1.70'; E='  +5.57%  '},
    @{Row=27; D='9.28'; E='  +7.79%  '},
    @{Row=28; E='  +2.53%  '},
    @{Row=29; D='8.09'; E='  +3.63%  '},
    @{Row=30; D='543.33'; E='  +1.53%  '},
    @{Row=31; E='  -0.01%  '},
    @{Row=32; E='  +2.45%  '},
    @{Row=33; E='  +7.67%  '},
    @{Row=34; E='  +1.06%  '},
    @{Row=35; D='5.26'; E='  +0.27%  '},
    @{Row=36; D='167.35'; E='  -0.78%  '},
    @{Row=37; D='2.01'; E='  +7.93%  '},
    @{Row=38; D='0.410'; E='  +2.36%  '},
    @{Row=39; D='1.00'; E='  +0.12%  '},
    @{Row=40; E='  +3.09%  '},
    @{Row=41; D='173.20'; E='  +3.39%  '},
    @{Row=42; E='  +0.05%  '},
    @{Row=43; D='40.05'; E='  +0.99%  '},
    @{Row=44; E='  +5.80%  '},
    @{Row=45; D='0.0588'; E='  +4.66%  '},
    @{Row=46; D='21.59'; E='  -2.51%  '},
    @{Row=47; D='0.630'; E='  +1.19%  '},
    @{Row=48; D='2.01'; E='  +15.62%  '},
    @{Row=49; D='0.0246'; E='  +2.63%  '},
    @{Row=50; D='0.0966'; E='  +1.25%  '},
    @{Row=51; D='19.33'; E='  +4.56%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('D')) {
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = '@'
        $cell.Value = $u.D
        $cell.Style = 'Normal'
    }
    if ($u.ContainsKey('E')) {
        $cell = $ws.Range("E$r")
        $cell.NumberFormat = '@'
        $cell.Value = $u.E
        $cell.Style = 'Normal'
    }
}
